$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number + week-covering dates) ---
$ws.Range("A8").Value = "Volume 32   Number  9"
$ws.Range("C9").Value = "Report Covering the Week  2/24/2025  Through  3/2/2025"

# --- Donor cells used to transplant exact style+content combos for cells that ---
# --- change between a plain number and the "N/A" shared-string placeholders  ---
# C14/D14/F14/I14 etc. hold style 13 + shared string "0"  (numeric placeholder text)
# E14/L14/M14 etc. hold style 13 + shared string "***.*" (non-numeric placeholder)
# G14 holds plain style 14 (integer format); H14 holds plain style 15 (decimal/pct format)
$naZero = $ws.Range("C14")     # style 13, text "0"
$naDash = $ws.Range("E14")     # style 13, text "***.*"
$numStyle = $ws.Range("G14")   # style 14 donor (integer cells)
$pctStyle = $ws.Range("H14")   # style 15 donor (percent/decimal cells)

# --- Cells changing type: number <-> "N/A" placeholder string ---
$naZero.Copy($ws.Range("F15"))
$naZero.Copy($ws.Range("C22"))
$numStyle.Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$pctStyle.Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$numStyle.Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1
$numStyle.Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 2
$pctStyle.Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = 0
$naZero.Copy($ws.Range("G29"))
$naDash.Copy($ws.Range("H29"))
$naZero.Copy($ws.Range("G30"))
$naDash.Copy($ws.Range("H30"))
$numStyle.Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D31").Value = 2
$pctStyle.Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E31").Value = -100
$numStyle.Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("G31").Value = 2
$pctStyle.Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("H31").Value = -100
$excel.CutCopyMode = $false

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("H15").Value = -100
$ws.Range("M15").Value = -80
$ws.Range("N15").Value = -93.333333333333
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 66.666666666666
$ws.Range("I16").Value = 42
$ws.Range("J16").Value = 34
$ws.Range("K16").Value = 23.529411764705
$ws.Range("L16").Value = -6.666666666666
$ws.Range("M16").Value = -12.5
$ws.Range("N16").Value = -82.857142857142
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -90
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = -51.724137931034
$ws.Range("I17").Value = 46
$ws.Range("J17").Value = 54
$ws.Range("K17").Value = -14.814814814814
$ws.Range("L17").Value = 2.222222222222
$ws.Range("M17").Value = 48.387096774193
$ws.Range("N17").Value = -60.344827586206
$ws.Range("C18").Value = 5
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 300
$ws.Range("I18").Value = 17
$ws.Range("K18").Value = 21.428571428571
$ws.Range("L18").Value = -26.086956521739
$ws.Range("M18").Value = -61.363636363636
$ws.Range("N18").Value = -94.848484848484
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = -23.529411764705
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = -38.59649122807
$ws.Range("I19").Value = 78
$ws.Range("J19").Value = 111
$ws.Range("K19").Value = -29.729729729729
$ws.Range("L19").Value = 1.298701298701
$ws.Range("M19").Value = 56
$ws.Range("N19").Value = -56.179775280898
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = -54.166666666666
$ws.Range("I20").Value = 33
$ws.Range("J20").Value = 39
$ws.Range("K20").Value = -15.384615384615
$ws.Range("L20").Value = -17.5
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = -90.406976744186
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -22.58064516129
$ws.Range("F21").Value = 92
$ws.Range("G21").Value = 127
$ws.Range("H21").Value = -27.55905511811
$ws.Range("I21").Value = 217
$ws.Range("J21").Value = 254
$ws.Range("K21").Value = -14.566929133858
$ws.Range("L21").Value = -6.465517241379
$ws.Range("M21").Value = 14.814814814814
$ws.Range("N21").Value = -82.542236524537
$ws.Range("J22").Value = 8
$ws.Range("K22").Value = -62.5
$ws.Range("L22").Value = -62.5
$ws.Range("M22").Value = -62.5
$ws.Range("D23").Value = 1
$ws.Range("F23").Value = 1
$ws.Range("H23").Value = -85.714285714285
$ws.Range("J23").Value = 9
$ws.Range("K23").Value = -66.666666666666
$ws.Range("M23").Value = -50
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -8.333333333333
$ws.Range("F24").Value = 79
$ws.Range("G24").Value = 97
$ws.Range("H24").Value = -18.556701030927
$ws.Range("I24").Value = 230
$ws.Range("J24").Value = 212
$ws.Range("K24").Value = 8.490566037735
$ws.Range("L24").Value = 9.523809523809
$ws.Range("M24").Value = 130
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 47
$ws.Range("H25").Value = -19.148936170212
$ws.Range("I25").Value = 88
$ws.Range("J25").Value = 93
$ws.Range("K25").Value = -5.376344086021
$ws.Range("L25").Value = -7.368421052631
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 225
$ws.Range("F26").Value = 44
$ws.Range("G26").Value = 44
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 88
$ws.Range("J26").Value = 83
$ws.Range("K26").Value = 6.024096385542
$ws.Range("L26").Value = 27.536231884058
$ws.Range("M26").Value = 2.325581395348
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 1
$ws.Range("I27").Value = 3
$ws.Range("K27").Value = 50
$ws.Range("L27").Value = -25
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 66.666666666666
$ws.Range("I28").Value = 7
$ws.Range("J28").Value = 8
$ws.Range("K28").Value = -12.5
$ws.Range("L28").Value = -12.5
$ws.Range("N29").Value = -97.5
$ws.Range("N30").Value = -97.435897435897
$ws.Range("J31").Value = 4
